$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1673003802281369
$ws.Range("C2").Value = 0.6463878326996197
$ws.Range("P2").Value = 0.1102661596958175
$ws.Range("S2").Value = 0.07604562737642585
$ws.Range("B3").Value = 0.005747126436781609
$ws.Range("C3").Value = 0.04597701149425287
$ws.Range("J3").Value = 0.005747126436781609
$ws.Range("P3").Value = 0.7586206896551724
$ws.Range("S3").Value = 0.1839080459770115
$ws.Range("P4").Value = 0.696969696969697
$ws.Range("S4").Value = 0.303030303030303
$ws.Range("B6").Value = 0.06493506493506493
$ws.Range("D6").Value = 0.008658008658008658
$ws.Range("F6").Value = 0.03463203463203463
$ws.Range("J6").Value = 0.2294372294372294
$ws.Range("O6").Value = 0.02597402597402598
$ws.Range("Q6").Value = 0.1818181818181818
$ws.Range("R6").Value = 0.05194805194805195
$ws.Range("S6").Value = 0.4025974025974026
$ws.Range("B7").Value = 0.1164021164021164
$ws.Range("D7").Value = 0.01587301587301587
$ws.Range("F7").Value = 0.04761904761904762
$ws.Range("J7").Value = 0.08994708994708994
$ws.Range("O7").Value = 0.01587301587301587
$ws.Range("Q7").Value = 0.253968253968254
$ws.Range("R7").Value = 0.1111111111111111
$ws.Range("S7").Value = 0.3492063492063492
$ws.Range("B8").Value = 0.08027522935779817
$ws.Range("D8").Value = 0.01146788990825688
$ws.Range("F8").Value = 0.07568807339449542
$ws.Range("J8").Value = 0.1032110091743119
$ws.Range("O8").Value = 0.009174311926605505
$ws.Range("Q8").Value = 0.1536697247706422
$ws.Range("R8").Value = 0.1009174311926606
$ws.Range("S8").Value = 0.4655963302752293
$ws.Range("B9").Value = 0.09352517985611511
$ws.Range("D9").Value = 0.02517985611510791
$ws.Range("E9").Value = 0.003597122302158274
$ws.Range("F9").Value = 0.07553956834532374
$ws.Range("J9").Value = 0.1043165467625899
$ws.Range("O9").Value = 0.02877697841726619
$ws.Range("Q9").Value = 0.1510791366906475
$ws.Range("R9").Value = 0.09352517985611511
$ws.Range("S9").Value = 0.4244604316546763
$ws.Range("B10").Value = 0.09882747068676717
$ws.Range("D10").Value = 0.01842546063651591
$ws.Range("F10").Value = 0.07537688442211055
$ws.Range("J10").Value = 0.09966499162479062
$ws.Range("O10").Value = 0.01340033500837521
$ws.Range("Q10").Value = 0.2185929648241206
$ws.Range("R10").Value = 0.0728643216080402
$ws.Range("S10").Value = 0.4028475711892797
$ws.Range("G11").Value = 0.1286173633440514
$ws.Range("J11").Value = 0.1254019292604502
$ws.Range("K11").Value = 0.1961414790996785
$ws.Range("L11").Value = 0.5337620578778135
$ws.Range("S11").Value = 0.01607717041800643
$ws.Range("G12").Value = 0.6983240223463687
$ws.Range("J12").Value = 0.223463687150838
$ws.Range("K12").Value = 0.01675977653631285
$ws.Range("L12").Value = 0.0335195530726257
$ws.Range("S12").Value = 0.02793296089385475
$ws.Range("G13").Value = 0.7631578947368421
$ws.Range("J13").Value = 0.2368421052631579
$ws.Range("F15").Value = 0.02953586497890295
$ws.Range("H15").Value = 0.1729957805907173
$ws.Range("I15").Value = 0.1308016877637131
$ws.Range("J15").Value = 0.3080168776371308
$ws.Range("K15").Value = 0.07172995780590717
$ws.Range("M15").Value = 0.008438818565400843
$ws.Range("O15").Value = 0.04219409282700422
$ws.Range("S15").Value = 0.2362869198312236
$ws.Range("F16").Value = 0.0111731843575419
$ws.Range("H16").Value = 0.1229050279329609
$ws.Range("I16").Value = 0.111731843575419
$ws.Range("J16").Value = 0.4022346368715084
$ws.Range("K16").Value = 0.09497206703910614
$ws.Range("M16").Value = 0.01675977653631285
$ws.Range("O16").Value = 0.0782122905027933
$ws.Range("S16").Value = 0.1620111731843575
$ws.Range("F17").Value = 0.02857142857142857
$ws.Range("H17").Value = 0.1626373626373626
$ws.Range("I17").Value = 0.1340659340659341
$ws.Range("J17").Value = 0.3802197802197802
$ws.Range("K17").Value = 0.0945054945054945
$ws.Range("M17").Value = 0.01978021978021978
$ws.Range("O17").Value = 0.07252747252747253
$ws.Range("S17").Value = 0.1076923076923077
$ws.Range("F18").Value = 0.0160427807486631
$ws.Range("H18").Value = 0.1818181818181818
$ws.Range("I18").Value = 0.09090909090909091
$ws.Range("J18").Value = 0.4010695187165775
$ws.Range("K18").Value = 0.0855614973262032
$ws.Range("M18").Value = 0.0106951871657754
$ws.Range("O18").Value = 0.05882352941176471
$ws.Range("S18").Value = 0.1550802139037433
$ws.Range("F19").Value = 0.01680058436815194
$ws.Range("H19").Value = 0.1928414901387874
$ws.Range("I19").Value = 0.1081081081081081
$ws.Range("J19").Value = 0.3389335281227173
$ws.Range("K19").Value = 0.1110299488677867
$ws.Range("M19").Value = 0.01680058436815194
$ws.Range("N19").Value = 0.002191380569758948
$ws.Range("O19").Value = 0.07523739956172389
$ws.Range("S19").Value = 0.1380569758948137
